# Updated cryptos list on Wed Sep 18 23:10:33 UTC 2024 with GitHub Actions
#
# Refreshes the Price / Volume(1h) columns with the latest scraped values,
# and reflects the new coinranking.com rank ordering where ImmutableX moved
# ahead of PolygonEcosystemToken (rows 34 and 35 swap places).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 34 & 35 swapped position (ImmutableX now ranks above PolygonEcosystemToken) ---
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "1.38"
$ws.Range("E34").Value = "  -0.42%  "

$ws.Range("B35").Value = "PolygonEcosystemToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D35").Value = "0.386"
$ws.Range("E35").Value = "  +1.41%  "

# --- Remaining Price / Volume(1h) refreshes ---
$ws.Range("D2").Value = "60.812.37"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "2.345.26"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "554.57"
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("D6").Value = "131.86"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "0.582"
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("D9").Value = "2.343.72"
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("D11").Value = "5.62"
$ws.Range("E11").Value = "  +2.15%  "
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").Value = "0.339"
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("D14").Value = "24.07"
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("D15").Value = "2.762.99"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").Value = "60.710.76"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("E17").Value = "  +1.85%  "
$ws.Range("D18").Value = "2.337.52"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "10.73"
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("D20").Value = "4.12"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "315.60"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").Value = "6.66"
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").Value = "64.32"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").Value = "8.01"
$ws.Range("E27").Value = "  +1.68%  "
$ws.Range("E28").Value = "  +4.78%  "
$ws.Range("D29").Value = "1.27"
$ws.Range("E29").Value = "  +8.83%  "
$ws.Range("D30").Value = "1.75"
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").Value = "0.0₃0739"
$ws.Range("E32").Value = "  +1.27%  "
$ws.Range("D33").Value = "6.14"
$ws.Range("E33").Value = "  +3.48%  "
$ws.Range("D36").Value = "18.10"
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("D39").Value = "4.16"
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("D40").Value = "334.73"
$ws.Range("E40").Value = "  +4.83%  "
$ws.Range("E41").Value = "  +1.48%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "139.82"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").Value = "3.53"
$ws.Range("E44").Value = "  +2.31%  "
$ws.Range("D45").Value = "0.0952"
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("D46").Value = "19.44"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "0.570"
$ws.Range("E47").Value = "  +2.18%  "
$ws.Range("D48").Value = "0.0501"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").Value = "0.0₆0228"
$ws.Range("E49").Value = "  +6.99%  "
$ws.Range("D50").Value = "0.0217"
$ws.Range("E50").Value = "  +2.41%  "
$ws.Range("D51").Value = "17.24"
$ws.Range("E51").Value = "  +2.50%  "
